# feat: add 2022-Q3 data
#
# The workbook tracks quarterly fund-holdings snapshots, one worksheet per
# quarter, plus a "总计" (summary) sheet. This edit inserts a brand-new
# "2022-Q3" quarter:
#   - a new row is inserted at the top of the "总计" summary table
#   - a new "2022-Q3" worksheet is inserted right after "总计" (i.e. right
#     before the existing "2022-Q2" sheet), built from a copy of the
#     "2022-Q2" sheet's layout/formatting, populated with the new quarter's
#     fund numbers.
# All the other quarterly sheets keep their own names/content unchanged;
# they simply shift one tab to the right to make room.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" summary sheet: insert a new second row for 2022-Q3.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Push the existing data rows down one, then restore the formatting
# (border/centering) that belongs to the index column A on the new row.
$summary.Rows.Item(2).Insert()
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.47

# Column A is a plain 0-based row index; fix it up across the whole table.
$summary.Range("A2").Value = 0
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6
$summary.Range("A9").Value = 7

# ---------------------------------------------------------------------
# 2) New "2022-Q3" worksheet, cloned from "2022-Q2" so it keeps the same
#    column layout/formatting, then updated with the new quarter's data.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# Row 2 - 513090 易方达中证香港证券投资主题ETF (name unchanged)
$q3.Range("D2").Value = "'10.53"
$q3.Range("E2").Value = "'96.33"
$q3.Range("F2").Value = "'3.72"
$q3.Range("G2").Value = "'0.3917"
$q3.Range("H2").Value = 10

# Row 3 - 011355 华泰柏瑞港股通时代机遇混合A
$q3.Range("C3").Value = "华泰柏瑞港股通时代机遇混合A"
$q3.Range("D3").Value = "'0.54"
$q3.Range("E3").Value = "'91.80"
$q3.Range("F3").Value = "'7.01"
$q3.Range("G3").Value = "'0.0379"
$q3.Range("H3").Value = 5

# Row 4 - 003413 华泰柏瑞新经济沪港深混合
$q3.Range("C4").Value = "华泰柏瑞新经济沪港深混合"
$q3.Range("D4").Value = "'0.42"
$q3.Range("E4").Value = "'86.45"
$q3.Range("F4").Value = "'6.53"
$q3.Range("G4").Value = "'0.0274"
$q3.Range("H4").Value = 1

# Row 5 - 011356 华泰柏瑞港股通时代机遇混合C
$q3.Range("C5").Value = "华泰柏瑞港股通时代机遇混合C"
$q3.Range("D5").Value = "'0.24"
$q3.Range("E5").Value = "'91.80"
$q3.Range("F5").Value = "'7.01"
$q3.Range("G5").Value = "'0.0168"
$q3.Range("H5").Value = 5
